$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 3, shifting existing data (old rows 3-13) down to rows 5-15.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# New row 3: Bing / Primera
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44901
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100103
$ws.Range("H3").Value = "Frutos de hueso (carozo)"
$ws.Range("I3").Value = 100103001
$ws.Range("J3").Value = "Cereza"
$ws.Range("K3").Value = "Bing"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 500
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 12500
$ws.Range("Q3").Value = "`$/caja 15 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 833
$ws.Range("T3").Value = 15

# New row 4: Lapins / Primera
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44901
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103001
$ws.Range("J4").Value = "Cereza"
$ws.Range("K4").Value = "Lapins"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 500
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 12500
$ws.Range("Q4").Value = "`$/caja 15 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 833
$ws.Range("T4").Value = 15
